$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the sheet from "Sheet1" to "Table" (this also updates the
# _xlnm._FilterDatabase defined name reference automatically).
$ws.Name = "Table"

# Update the active cell selection on the sheet.
$ws.Range("E12").Select()
